$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("demand_nodes")

# Update the max_release (E) column values for the demand nodes
$ws.Range("E2").Value = 22
$ws.Range("E3").Value = 28
$ws.Range("E4").Value = 27
$ws.Range("E5").Value = 22
$ws.Range("E6").Value = 19
$ws.Range("E7").Value = 37
$ws.Range("E8").Value = 39
$ws.Range("E9").Value = 25

# Apply black font color to a new block of cells below the table (E14:E21)
$range = $ws.Range("E14:E21")
$range.Font.Color = 0

# Update the active selection to mirror the author's final cursor position
$ws.Range("C16").Select()
